$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 115, pushing the existing data (rows 115-252) down by one.
$ws.Rows.Item(115).Insert()

# The record that used to live on row 115 is now on row 116; duplicate it back onto the
# newly inserted row 115 so the new row starts out identical to the record above it.
$ws.Range("A116:R116").Copy()
$ws.Range("A115").PasteSpecial()

# Row 115 represents the newest weekly observation, so give it its own (later) date while
# keeping every other field the same as the record it was copied from.
$ws.Range("D115").Value = 44671
